$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.885.65'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '2.325.73'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.17'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.43%  '
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').Value = '2.690.11'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = '2.326.60'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.791'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('D18').Value = '42.818.84'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.62%  '
$ws.Range('D21').Value = '0.0₃0893'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').Value = '  +3.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('E38').Value = '  +3.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.1000'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +13.06%  '
$ws.Range('D43').Value = '1.943.52'
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.09'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.04%  '
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('D48').Value = '2.556.66'
$ws.Range('E48').Value = '  +1.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.53'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.71%  '
